$wb = $excel.ActiveWorkbook

# --- 1. Update Metadata sheet values (Version, Date, Contact) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-03T10:45:43+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- 2. Update the C2 "code" cell on the existing 5 "Include from FSIII" sheets ---
$codeMap = @{
    "Include from FSIII"   = "43c2b7f0-5e55-4627-8fcf-bdaf5a9d84ac"
    "Include from FSIII 2" = "1c850a09-aa49-4fae-9354-f932f13e030b"
    "Include from FSIII 3" = "462f9352-0129-4d8e-8c75-a6dfed78ddcf"
    "Include from FSIII 4" = "4571f168-a92a-4caf-8dc8-35f45c2a1cb4"
    "Include from FSIII 5" = "86b53158-6d05-412e-ad55-2e1fa26359b3"
}

foreach ($name in @("Include from FSIII","Include from FSIII 2","Include from FSIII 3","Include from FSIII 4","Include from FSIII 5")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Range("C2").Value = $codeMap[$name]
}

# --- 3. Add 5 new sheets (copies of the existing layout), with the J1..J5 codes ---
$newCodes = @{
    "Include from FSIII 6"  = "J1"
    "Include from FSIII 7"  = "J2"
    "Include from FSIII 8"  = "J3"
    "Include from FSIII 9"  = "J4"
    "Include from FSIII 10" = "J5"
}

foreach ($name in @("Include from FSIII 6","Include from FSIII 7","Include from FSIII 8","Include from FSIII 9","Include from FSIII 10")) {
    $template = $wb.Worksheets.Item("Include from FSIII 5")
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $last)
    $new = $wb.Worksheets.Item($wb.Worksheets.Count)
    $new.Name = $name
    $new.Range("C2").Value = $newCodes[$name]
}
